# Auto-generated from the XML diff: updates cell text values on the
# active sheet of the cryptos workbook. A leading apostrophe forces
# Excel to store the value as text (matching the original t="inlineStr"
# / shared-string cell type) instead of auto-coercing numeric-looking
# strings (e.g. "521.18") into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.381.14"
$ws.Range("E2").Value = "'  +0.89%  "
$ws.Range("D3").Value = "'2.480.57"
$ws.Range("E3").Value = "'  +1.11%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'521.18"
$ws.Range("E5").Value = "'  +0.77%  "
$ws.Range("D6").Value = "'132.62"
$ws.Range("E6").Value = "'  +0.57%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "'  -0.23%  "
$ws.Range("D8").Value = "'0.556"
$ws.Range("E8").Value = "'  +0.15%  "
$ws.Range("D9").Value = "'2.511.16"
$ws.Range("E9").Value = "'  +2.18%  "
$ws.Range("D10").Value = "'0.0976"
$ws.Range("E10").Value = "'  -0.47%  "
$ws.Range("E11").Value = "'  -0.24%  "
$ws.Range("E12").Value = "'  -1.87%  "
$ws.Range("E13").Value = "'  -2.20%  "
$ws.Range("D14").Value = "'2.925.98"
$ws.Range("E14").Value = "'  +1.25%  "
$ws.Range("D15").Value = "'58.316.78"
$ws.Range("E15").Value = "'  +0.91%  "
$ws.Range("D16").Value = "'22.14"
$ws.Range("E16").Value = "'  -0.21%  "
$ws.Range("E17").Value = "'  +0.00%  "
$ws.Range("D18").Value = "'2.501.41"
$ws.Range("E18").Value = "'  +1.79%  "
$ws.Range("D19").Value = "'10.71"
$ws.Range("E19").Value = "'  +0.64%  "
$ws.Range("D20").Value = "'320.93"
$ws.Range("E20").Value = "'  +0.57%  "
$ws.Range("E21").Value = "'  +0.47%  "
$ws.Range("D22").Value = "'6.02"
$ws.Range("E22").Value = "'  +5.23%  "
$ws.Range("E23").Value = "'  -0.38%  "
$ws.Range("D24").Value = "'63.86"
$ws.Range("E24").Value = "'  -0.66%  "
$ws.Range("D25").Value = "'0.400"
$ws.Range("E25").Value = "'  -1.80%  "
$ws.Range("E26").Value = "'  +0.93%  "
$ws.Range("D27").Value = "'0.991"
$ws.Range("E27").Value = "'  -0.67%  "
$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = "'  +0.76%  "
$ws.Range("D29").Value = "'0.0₃0755"
$ws.Range("E29").Value = "'  +2.68%  "
$ws.Range("B30").Value = "'PancakeSwap"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.71"
$ws.Range("E30").Value = "'  +1.49%  "
$ws.Range("B31").Value = "'Fetch.AI"
$ws.Range("C31").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.20"
$ws.Range("E31").Value = "'  +2.99%  "
$ws.Range("D32").Value = "'167.21"
$ws.Range("E32").Value = "'  +0.90%  "
$ws.Range("D33").Value = "'6.25"
$ws.Range("E33").Value = "'  +0.76%  "
$ws.Range("D35").Value = "'0.992"
$ws.Range("E35").Value = "'  -0.49%  "
$ws.Range("D36").Value = "'18.05"
$ws.Range("E36").Value = "'  +0.17%  "
$ws.Range("D37").Value = "'1.26"
$ws.Range("E37").Value = "'  -2.37%  "
$ws.Range("E38").Value = "'  +0.10%  "
$ws.Range("D39").Value = "'36.87"
$ws.Range("E39").Value = "'  +1.84%  "
$ws.Range("D40").Value = "'1.47"
$ws.Range("E40").Value = "'  -0.46%  "
$ws.Range("D41").Value = "'0.778"
$ws.Range("E41").Value = "'  -0.86%  "
$ws.Range("D42").Value = "'277.74"
$ws.Range("E42").Value = "'  +2.59%  "
$ws.Range("D43").Value = "'5.08"
$ws.Range("E43").Value = "'  +1.53%  "
$ws.Range("D44").Value = "'3.43"
$ws.Range("E44").Value = "'  +0.22%  "
$ws.Range("D45").Value = "'0.597"
$ws.Range("E45").Value = "'  +1.57%  "
$ws.Range("D46").Value = "'122.33"
$ws.Range("E46").Value = "'  -1.74%  "
$ws.Range("D47").Value = "'0.0917"
$ws.Range("E47").Value = "'  +1.36%  "
$ws.Range("D48").Value = "'0.0501"
$ws.Range("E48").Value = "'  +3.29%  "
$ws.Range("D49").Value = "'17.81"
$ws.Range("E49").Value = "'  +1.12%  "
$ws.Range("E50").Value = "'  +1.65%  "
$ws.Range("D51").Value = "'16.91"
$ws.Range("E51").Value = "'  +1.45%  "
